$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the long preparation text in F2 to reflect the new transaction code (090 -> 998)
$ws.Range("F2").Value = "Username : 30711;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nKode Transaksi : 998;`nNama Jenis Transaksi : Saldo Awal Top Up;`nGroup Transaksi : IBA;`nKontributor : EE;`nKode Dokumen : 002 : ID Peserta"

# Update the "Kode Transaksi" value from text "090" to the number 998, while keeping
# the cell's existing number format/style (it used a quote-prefixed style since "090"
# looked numeric). Setting .Value directly resets the style, so write the value then
# restore the original format via a format-only paste from a cell sharing that style.
$ws.Range("M2").Value = 998
$ws.Range("K2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the view: scroll position and selection (was I1/Q1:R1, now E1/G2)
$ws.Range("E1").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("G2").Select() | Out-Null
